# Auto-update: decrement column E (剩余/"remaining") by 1 for every data row
# (rows 2-99), except row 36 whose value is left untouched, matching the
# commit "自动更新Excel文件 - 2025-10-14 23:11:51".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }

    $cell = $ws.Cells.Item($row, 5)
    $current = $cell.Value2

    if ($null -ne $current) {
        $cell.Value = $current - 1
    }
}
